# Generate Report for Handback
#
# This mirrors a "handback" event being recorded in the localization-status
# workbook: the two language sheets (zh-cn, de-de) each get their first data
# row (the 5e667d6f-... file) filled in with "Latest Target File" / "Latest
# Handback File" hyperlinks plus an updated "Latest Handback DateTime", and
# the shared "Status" text flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview +
# both language sheets).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile  = "5e667d6f-dc90-457e-b295-fbe6abeb0028.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/5e667d6f-dc90-457e-b295-fbe6abeb0028.md"

$zhFile  = "5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf"
$zhUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b73bc382c3f9b832b82cac15f26298cfcaba2b92/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf"

$deFile  = "5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf"
$deUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25400aeee798a320ea462dfbcc625c51a5a62fba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf"

# ---- 1. Flip the shared "Status" text everywhere it shows up ----

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $statusNew
$zh.Range("B3").Value = $statusNew

$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $statusNew
$de.Range("B3").Value = $statusNew

# ---- 2. zh-cn sheet: fill in E2/F2 (Target/Handback files) + G2 datetime ----

$zh.Range("E2").Value = $mdFile
$zh.Range("E2").Hyperlinks.Add($zh.Range("E2"), $mdUrl, "", "", $mdFile) | Out-Null
$zh.Range("E2").Style = "HyperLink"

$zh.Range("F2").Value = $zhFile
$zh.Range("F2").Hyperlinks.Add($zh.Range("F2"), $zhUrl, "", "", $zhFile) | Out-Null
$zh.Range("F2").Style = "HyperLink"

$zh.Range("G2").Value = "2016-01-28 05:52:30"

$zh.Range("E3").Value = $mdFile
$zh.Range("E3").Hyperlinks.Add($zh.Range("E3"), $mdUrl, "", "", $mdFile) | Out-Null
$zh.Range("E3").Style = "HyperLink"

$zh.Range("F3").Value = $zhFile
$zh.Range("F3").Hyperlinks.Add($zh.Range("F3"), $zhUrl, "", "", $zhFile) | Out-Null
$zh.Range("F3").Style = "HyperLink"

# ---- 3. de-de sheet: fill in E2/F2 (Target/Handback files) + G2 datetime ----

$de.Range("E2").Value = $mdFile
$de.Range("E2").Hyperlinks.Add($de.Range("E2"), $mdUrl, "", "", $mdFile) | Out-Null
$de.Range("E2").Style = "HyperLink"

$de.Range("F2").Value = $deFile
$de.Range("F2").Hyperlinks.Add($de.Range("F2"), $deUrl, "", "", $deFile) | Out-Null
$de.Range("F2").Style = "HyperLink"

$de.Range("G2").Value = "2016-01-28 05:52:47"

$de.Range("E3").Value = $mdFile
$de.Range("E3").Hyperlinks.Add($de.Range("E3"), $mdUrl, "", "", $mdFile) | Out-Null
$de.Range("E3").Style = "HyperLink"

$de.Range("F3").Value = $deFile
$de.Range("F3").Hyperlinks.Add($de.Range("F3"), $deUrl, "", "", $deFile) | Out-Null
$de.Range("F3").Style = "HyperLink"
